# Apply the edits described by the diff:
# 1. Update the confidential disclaimer text date from 2021-04-27 to 2021-04-28
# 2. Update numeric values in columns D and E for rows 2-38

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Unprotect()

# --- 1. Update the disclaimer text (cell A41) ---
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."
$ws.Range("A41").Value = $newText

# --- 2. Update numeric D/E columns for rows 2-38 ---
$ws.Range("D2").Value = 0.03012474032144049
$ws.Range("E2").Value = 0.00439698492462326
$ws.Range("D3").Value = 0.02996640856659752
$ws.Range("E3").Value = -0.003556792700842748
$ws.Range("D4").Value = 0.03126067411868587
$ws.Range("E4").Value = -0.01396549700739347
$ws.Range("D5").Value = 0.06598630353085606
$ws.Range("E5").Value = 0.01201780285185072
$ws.Range("D6").Value = 0.014997106583724
$ws.Range("E6").Value = 0.001467748165314697
$ws.Range("D7").Value = 0.01598880401405754
$ws.Range("E7").Value = -0.0104340265198174
$ws.Range("D8").Value = 0.02955744922908848
$ws.Range("E8").Value = -0.07219848704581977
$ws.Range("D9").Value = 0.03421124429643644
$ws.Range("E9").Value = 0.00496670053053383
$ws.Range("D10").Value = 0.02928461658324321
$ws.Range("E10").Value = 0.001918702403322925
$ws.Range("D11").Value = 0.03188724307535105
$ws.Range("E11").Value = 0.003390979993218135
$ws.Range("D12").Value = 0.01368700785615135
$ws.Range("E12").Value = -0.01918600550186922
$ws.Range("D13").Value = 0.01430276391248579
$ws.Range("E13").Value = -0.01606501606501598
$ws.Range("D14").Value = 0.01631396337400335
$ws.Range("E14").Value = -0.001988400994200612
$ws.Range("D15").Value = 0.008322843854574987
$ws.Range("E15").Value = 0.03489235337787666
$ws.Range("D16").Value = 0.007078008703998591
$ws.Range("E16").Value = 0.05646943285047867
$ws.Range("D17").Value = 0.03165611733078149
$ws.Range("E17").Value = 0.01140002561803488
$ws.Range("D18").Value = 0.02980769063674274
$ws.Range("E18").Value = -0.01303976058144518
$ws.Range("D19").Value = 0.03232864111385231
$ws.Range("E19").Value = -0.02365167532700241
$ws.Range("D20").Value = 0.02930778708395194
$ws.Range("E20").Value = 0.01162829001548249
$ws.Range("D21").Value = 0.04454760466259985
$ws.Range("E21").Value = 0.03155015777246084
$ws.Range("D22").Value = 0.0334649610861095
$ws.Range("E22").Value = 0.004269682370250827
$ws.Range("D23").Value = 0.03095096175921254
$ws.Range("E23").Value = -0.002744939018684223
$ws.Range("D24").Value = 0.02947982805171424
$ws.Range("E24").Value = -0.004244282008959988
$ws.Range("D25").Value = 0.01522823232829356
$ws.Range("E25").Value = -0.001825858724181173
$ws.Range("D26").Value = 0.01480015732769981
$ws.Range("E26").Value = 0.01972602739726037
$ws.Range("D27").Value = 0.03016239238509218
$ws.Range("E27").Value = 0.006971340046475705
$ws.Range("D28").Value = 0.02924580599455609
$ws.Range("E28").Value = -0.007341678550678798
$ws.Range("D29").Value = 0.03034988035332697
$ws.Range("E29").Value = -0.02828568156659161
$ws.Range("D30").Value = 0.02805966944577514
$ws.Range("E30").Value = -0.0105972295814093
$ws.Range("D31").Value = 0.03564028492764722
$ws.Range("E31").Value = -0.006826271393046923
$ws.Range("D32").Value = 0.03034930109080925
$ws.Range("E32").Value = -0.007418293792427644
$ws.Range("D33").Value = 0.03044043839359691
$ws.Range("E33").Value = -0.01327616063329762
$ws.Range("D34").Value = 0.03046920843197691
$ws.Range("E34").Value = 0.0002281368821293039
$ws.Range("D35").Value = 0.02996988414170383
$ws.Range("E35").Value = 0.0002319378406585582
$ws.Range("D36").Value = 0.02920448526829219
$ws.Range("E36").Value = 0.006479338842975135
$ws.Range("D37").Value = 0.03156749016557061
$ws.Range("E37").Value = 0.0004893325504011781
$ws.Range("E38").Value = -0.002060436775523478

# Restore sheet protection (the sheet was protected before the edit; the
# diff does not indicate any intentional change to protection state, so we
# re-lock it after writing the new values). The original legacy password
# hash cannot be reproduced bit-for-bit through COM automation, so we
# simply restore the protected state.
$ws.Protect()
